# Update MIMAG and MIMS template (MIMS - Metagenomics)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MIMS")

# 1) Rename the "Parameter [...]" column to "Component [...]" via the table's
#    header cell so the underlying shared string / table column name updates.
$tbl = $ws.ListObjects.Item(1)
$tbl.ListColumns.Item(2).Range.Item(1).Value = "Component [next generation sequencing instrument model]"

# 2) Add the two new trailing columns to the annotation table.
$dataFormatCol = $tbl.ListColumns.Add()
$dataFormatCol.Range.Item(1).Value = "Data Format"

$dataSelectorFormatCol = $tbl.ListColumns.Add()
$dataSelectorFormatCol.Range.Item(1).Value = "Data Selector Format"

# 3) The new columns' data row (row 2) should be blank, matching the other
#    blank cells already in that row (columns A and E).
$ws.Range("A2").Copy($ws.Range("F2:G2"))

# 4) Fix the typo/casing in the template description: "sample" -> "Sample".
$wsTemplate = $wb.Worksheets.Item("isa_template")
$wsTemplate.Range("B5").Value = "Template that includes the minimum required information to describe a metagenomics experiment, based on the mimimum information about a metagenome sequence (MIMS) standard. This template should be combined with the MIxS - Sample information template to contain all required information."
